# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.968.30"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "'2.931.51"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'372.96"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'100.35"
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("D7").Value = "'0.533"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "'35.99"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "'0.138"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "'0.0844"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'3.399.54"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'17.88"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "'7.47"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "'2.936.17"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'10.98"
$ws.Range("E17").Value = "  +47.69%  "
$ws.Range("D18").Value = "'0.966"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'50.971.59"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "'3.14"
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("D21").Value = "'12.32"
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("D22").Value = "'0.0₃0953"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'68.55"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'263.96"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'3.10"
$ws.Range("E25").Value = "  +9.70%  "
$ws.Range("D26").Value = "'8.03"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'25.53"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.110"
$ws.Range("E30").Value = "  -4.36%  "
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").Value = "'9.89"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'50.65"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").Value = "'33.05"
$ws.Range("D36").Value = "'0.0440"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'16.24"
$ws.Range("E40").Value = "  -5.40%  "
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("D42").Value = "'2.46"
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").Value = "'119.68"
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").Value = "'21.03"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.03"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.271"
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("D47").Value = "'3.28"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "'1.984.18"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").Value = "'0.0325"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").Value = "'1.29"
$ws.Range("E51").Value = "  +0.27%  "
